$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1675
$ws.Range("F3").Value = 9581
$ws.Range("F4").Value = 781
$ws.Range("F5").Value = 632
$ws.Range("F6").Value = 248
$ws.Range("F7").Value = 330
$ws.Range("F10").Value = 1442
$ws.Range("F11").Value = 579
$ws.Range("F12").Value = 63
$ws.Range("F13").Value = 1509
$ws.Range("F14").Value = 132
$ws.Range("F15").Value = 314
$ws.Range("F17").Value = 157
$ws.Range("F18").Value = 93
$ws.Range("F19").Value = 415
$ws.Range("F20").Value = 1116
$ws.Range("F21").Value = 107
$ws.Range("F25").Value = 292
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 272
$ws.Range("F29").Value = 610
$ws.Range("F31").Value = 10
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 183
$ws.Range("F35").Value = 34
$ws.Range("F36").Value = 190
$ws.Range("F37").Value = 339
$ws.Range("F38").Value = 514
$ws.Range("F39").Value = 332
$ws.Range("F40").Value = 636
$ws.Range("F43").Value = 333
$ws.Range("F44").Value = 286
$ws.Range("F45").Value = 331
$ws.Range("F46").Value = 62
$ws.Range("F47").Value = 330

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 17
$ws.Range("F11").Value = 703
$ws.Range("F19").Value = 983
$ws.Range("F21").Value = 1072
$ws.Range("F22").Value = 296
$ws.Range("F23").Value = 662
$ws.Range("F24").Value = 33
$ws.Range("F27").Value = 333
$ws.Range("F32").Value = 127
$ws.Range("F35").Value = 24
$ws.Range("F37").Value = 118
$ws.Range("F39").Value = 33
$ws.Range("F40").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 781
$ws.Range("F7").Value = 2423
$ws.Range("F8").Value = 3751
$ws.Range("F11").Value = 136
$ws.Range("F12").Value = 122

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1675
$ws.Range("F3").Value = 17
$ws.Range("F4").Value = 781
$ws.Range("F5").Value = 9581
$ws.Range("F6").Value = 3751
$ws.Range("F7").Value = 781
$ws.Range("F8").Value = 136
$ws.Range("F9").Value = 136
$ws.Range("F11").Value = 632
$ws.Range("F12").Value = 248
$ws.Range("F13").Value = 330
$ws.Range("F15").Value = 703
$ws.Range("F16").Value = 1442
$ws.Range("F17").Value = 579
$ws.Range("F18").Value = 122
$ws.Range("F19").Value = 122
$ws.Range("F20").Value = 1509
$ws.Range("F21").Value = 132
$ws.Range("F22").Value = 314
$ws.Range("F24").Value = 157
$ws.Range("F25").Value = 1116
$ws.Range("F26").Value = 107
$ws.Range("F29").Value = 292
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 272
$ws.Range("F33").Value = 1072
$ws.Range("F34").Value = 296
$ws.Range("F35").Value = 610
$ws.Range("F37").Value = 33
$ws.Range("F39").Value = 333
$ws.Range("F40").Value = 339
$ws.Range("F41").Value = 514
$ws.Range("F42").Value = 332
$ws.Range("F44").Value = 636
$ws.Range("F47").Value = 333
$ws.Range("F48").Value = 118
$ws.Range("F49").Value = 331
$ws.Range("F50").Value = 330
$ws.Range("F51").Value = 33
